$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 131 (a new week's price
# observations), pushing the existing rows 131-134 down to 133-136.
$ws.Rows.Item(131).Insert()
$ws.Rows.Item(131).Insert()

# The newly inserted rows are blank; seed them from the rows that used to
# be directly above them (now at 133/134) so every column - including
# styling such as the date format on column D - matches the existing data
# pattern, then overwrite just the cells that actually carry new figures.
$ws.Range("A133:T133").Copy()
$ws.Range("A131:T131").PasteSpecial()
$ws.Range("A134:T134").Copy()
$ws.Range("A132:T132").PasteSpecial()

# Row 131 (Primera): new date + updated volume/price figures.
$ws.Range("D131").Value = 44448
$ws.Range("N131").Value = 12000
$ws.Range("O131").Value = 12500
$ws.Range("P131").Value = 12250
$ws.Range("S131").Value = 681

# Row 132 (Segunda): new date + updated volume/price figures.
$ws.Range("D132").Value = 44448
$ws.Range("M132").Value = 60
$ws.Range("N132").Value = 11000
$ws.Range("O132").Value = 11500
$ws.Range("P132").Value = 11250
$ws.Range("S132").Value = 625
